$d = $word.ActiveDocument

$d.Content.Find.Execute("82+8=", $true, $false, $false, $false, $false, $true, 1, $false, "98-50=", 2) | Out-Null
$d.Content.Find.Execute("16-12=", $true, $false, $false, $false, $false, $true, 1, $false, "58+24=", 2) | Out-Null
$d.Content.Find.Execute("36+41=", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=", 2) | Out-Null
$d.Content.Find.Execute("98-52=", $true, $false, $false, $false, $false, $true, 1, $false, "18+29=", 2) | Out-Null
$d.Content.Find.Execute("60-36=", $true, $false, $false, $false, $false, $true, 1, $false, "67-24=", 2) | Out-Null
$d.Content.Find.Execute("61+10=", $true, $false, $false, $false, $false, $true, 1, $false, "42+8=", 2) | Out-Null
$d.Content.Find.Execute("44+53=", $true, $false, $false, $false, $false, $true, 1, $false, "23+34=", 2) | Out-Null
$d.Content.Find.Execute("85-8=", $true, $false, $false, $false, $false, $true, 1, $false, "17+41=", 2) | Out-Null
$d.Content.Find.Execute("28-21=", $true, $false, $false, $false, $false, $true, 1, $false, "92-9=", 2) | Out-Null
$d.Content.Find.Execute("8-1=", $true, $false, $false, $false, $false, $true, 1, $false, "30+63=", 2) | Out-Null
$d.Content.Find.Execute("66-33=", $true, $false, $false, $false, $false, $true, 1, $false, "91-64=", 2) | Out-Null
$d.Content.Find.Execute("10-4=", $true, $false, $false, $false, $false, $true, 1, $false, "0+97=", 2) | Out-Null
$d.Content.Find.Execute("30+51=", $true, $false, $false, $false, $false, $true, 1, $false, "25-6=", 2) | Out-Null
$d.Content.Find.Execute("74-47=", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=", 2) | Out-Null
$d.Content.Find.Execute("24+17=", $true, $false, $false, $false, $false, $true, 1, $false, "53-26=", 2) | Out-Null
$d.Content.Find.Execute("63-31=", $true, $false, $false, $false, $false, $true, 1, $false, "4+77=", 2) | Out-Null
$d.Content.Find.Execute("65-28=", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=", 2) | Out-Null
$d.Content.Find.Execute("89-27=", $true, $false, $false, $false, $false, $true, 1, $false, "11+23=", 2) | Out-Null
$d.Content.Find.Execute("87-22=", $true, $false, $false, $false, $false, $true, 1, $false, "53-40=", 2) | Out-Null
$d.Content.Find.Execute("44+26=", $true, $false, $false, $false, $false, $true, 1, $false, "44+39=", 2) | Out-Null
$d.Content.Find.Execute("67-17=", $true, $false, $false, $false, $false, $true, 1, $false, "19-14=", 2) | Out-Null
$d.Content.Find.Execute("81-50=", $true, $false, $false, $false, $false, $true, 1, $false, "31+24=", 2) | Out-Null
$d.Content.Find.Execute("67-28=", $true, $false, $false, $false, $false, $true, 1, $false, "66-20=", 2) | Out-Null
$d.Content.Find.Execute("65-24=", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=", 2) | Out-Null
$d.Content.Find.Execute("11+29=", $true, $false, $false, $false, $false, $true, 1, $false, "66+32=", 2) | Out-Null
$d.Content.Find.Execute("72-72=", $true, $false, $false, $false, $false, $true, 1, $false, "9+81=", 2) | Out-Null
$d.Content.Find.Execute("88-86=", $true, $false, $false, $false, $false, $true, 1, $false, "4+48=", 2) | Out-Null
$d.Content.Find.Execute("5+7=", $true, $false, $false, $false, $false, $true, 1, $false, "13+60=", 2) | Out-Null
$d.Content.Find.Execute("88-28=", $true, $false, $false, $false, $false, $true, 1, $false, "28+47=", 2) | Out-Null
$d.Content.Find.Execute("2+97=", $true, $false, $false, $false, $false, $true, 1, $false, "17+8=", 2) | Out-Null
$d.Content.Find.Execute("62+6=", $true, $false, $false, $false, $false, $true, 1, $false, "14-11=", 2) | Out-Null
$d.Content.Find.Execute("25+54=", $true, $false, $false, $false, $false, $true, 1, $false, "51-6=", 2) | Out-Null
$d.Content.Find.Execute("40-40=", $true, $false, $false, $false, $false, $true, 1, $false, "11+2=", 2) | Out-Null
$d.Content.Find.Execute("8+9=", $true, $false, $false, $false, $false, $true, 1, $false, "8+89=", 2) | Out-Null
$d.Content.Find.Execute("32+61=", $true, $false, $false, $false, $false, $true, 1, $false, "14+56=", 2) | Out-Null
$d.Content.Find.Execute("66-18=", $true, $false, $false, $false, $false, $true, 1, $false, "91+3=", 2) | Out-Null
$d.Content.Find.Execute("69+3=", $true, $false, $false, $false, $false, $true, 1, $false, "58-18=", 2) | Out-Null
$d.Content.Find.Execute("47-0=", $true, $false, $false, $false, $false, $true, 1, $false, "66-36=", 2) | Out-Null
$d.Content.Find.Execute("74-7=", $true, $false, $false, $false, $false, $true, 1, $false, "8+87=", 2) | Out-Null
$d.Content.Find.Execute("25+34=", $true, $false, $false, $false, $false, $true, 1, $false, "88-45=", 2) | Out-Null
$d.Content.Find.Execute("46-15=", $true, $false, $false, $false, $false, $true, 1, $false, "86-60=", 2) | Out-Null
$d.Content.Find.Execute("93-83=", $true, $false, $false, $false, $false, $true, 1, $false, "20+50=", 2) | Out-Null
$d.Content.Find.Execute("54-31=", $true, $false, $false, $false, $false, $true, 1, $false, "32+6=", 2) | Out-Null
$d.Content.Find.Execute("79-58=", $true, $false, $false, $false, $false, $true, 1, $false, "73+11=", 2) | Out-Null
$d.Content.Find.Execute("0+17=", $true, $false, $false, $false, $false, $true, 1, $false, "92-22=", 2) | Out-Null
$d.Content.Find.Execute("82-2=", $true, $false, $false, $false, $false, $true, 1, $false, "59-57=", 2) | Out-Null
$d.Content.Find.Execute("72-69=", $true, $false, $false, $false, $false, $true, 1, $false, "53-42=", 2) | Out-Null
$d.Content.Find.Execute("34+55=", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=", 2) | Out-Null
$d.Content.Find.Execute("36-22=", $true, $false, $false, $false, $false, $true, 1, $false, "2+69=", 2) | Out-Null
$d.Content.Find.Execute("22-21=", $true, $false, $false, $false, $false, $true, 1, $false, "18+42=", 2) | Out-Null
$d.Content.Find.Execute("7+4=", $true, $false, $false, $false, $false, $true, 1, $false, "88-30=", 2) | Out-Null
$d.Content.Find.Execute("1+39=", $true, $false, $false, $false, $false, $true, 1, $false, "50-13=", 2) | Out-Null
$d.Content.Find.Execute("3+12=", $true, $false, $false, $false, $false, $true, 1, $false, "85-83=", 2) | Out-Null
$d.Content.Find.Execute("88+3=", $true, $false, $false, $false, $false, $true, 1, $false, "36+11=", 2) | Out-Null
$d.Content.Find.Execute("9+63=", $true, $false, $false, $false, $false, $true, 1, $false, "53+18=", 2) | Out-Null
$d.Content.Find.Execute("28+26=", $true, $false, $false, $false, $false, $true, 1, $false, "23+8=", 2) | Out-Null
$d.Content.Find.Execute("68-42=", $true, $false, $false, $false, $false, $true, 1, $false, "61+38=", 2) | Out-Null
$d.Content.Find.Execute("14+61=", $true, $false, $false, $false, $false, $true, 1, $false, "47+5=", 2) | Out-Null
$d.Content.Find.Execute("21+65=", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=", 2) | Out-Null
$d.Content.Find.Execute("95-11=", $true, $false, $false, $false, $false, $true, 1, $false, "46-37=", 2) | Out-Null
$d.Content.Find.Execute("24+13=", $true, $false, $false, $false, $false, $true, 1, $false, "87-35=", 2) | Out-Null
$d.Content.Find.Execute("77-28=", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=", 2) | Out-Null
$d.Content.Find.Execute("92-68=", $true, $false, $false, $false, $false, $true, 1, $false, "29+64=", 2) | Out-Null
$d.Content.Find.Execute("22+11=", $true, $false, $false, $false, $false, $true, 1, $false, "17-11=", 2) | Out-Null
$d.Content.Find.Execute("36-13=", $true, $false, $false, $false, $false, $true, 1, $false, "10+21=", 2) | Out-Null
$d.Content.Find.Execute("41+35=", $true, $false, $false, $false, $false, $true, 1, $false, "15+80=", 2) | Out-Null
$d.Content.Find.Execute("51-30=", $true, $false, $false, $false, $false, $true, 1, $false, "76+0=", 2) | Out-Null
$d.Content.Find.Execute("37+24=", $true, $false, $false, $false, $false, $true, 1, $false, "54+10=", 2) | Out-Null
$d.Content.Find.Execute("57+9=", $true, $false, $false, $false, $false, $true, 1, $false, "33+13=", 2) | Out-Null
$d.Content.Find.Execute("95-64=", $true, $false, $false, $false, $false, $true, 1, $false, "86-55=", 2) | Out-Null
$d.Content.Find.Execute("81-35=", $true, $false, $false, $false, $false, $true, 1, $false, "12+86=", 2) | Out-Null
$d.Content.Find.Execute("16+46=", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=", 2) | Out-Null
$d.Content.Find.Execute("20+5=", $true, $false, $false, $false, $false, $true, 1, $false, "39+50=", 2) | Out-Null
$d.Content.Find.Execute("57+7=", $true, $false, $false, $false, $false, $true, 1, $false, "4-4=", 2) | Out-Null
$d.Content.Find.Execute("53-46=", $true, $false, $false, $false, $false, $true, 1, $false, "54-33=", 2) | Out-Null
$d.Content.Find.Execute("44-27=", $true, $false, $false, $false, $false, $true, 1, $false, "89-74=", 2) | Out-Null
$d.Content.Find.Execute("38+22=", $true, $false, $false, $false, $false, $true, 1, $false, "98-98=", 2) | Out-Null
$d.Content.Find.Execute("11-10=", $true, $false, $false, $false, $false, $true, 1, $false, "71+13=", 2) | Out-Null
$d.Content.Find.Execute("3+67=", $true, $false, $false, $false, $false, $true, 1, $false, "77+7=", 2) | Out-Null
$d.Content.Find.Execute("59-45=", $true, $false, $false, $false, $false, $true, 1, $false, "25-17=", 2) | Out-Null
$d.Content.Find.Execute("32+18=", $true, $false, $false, $false, $false, $true, 1, $false, "71-21=", 2) | Out-Null
$d.Content.Find.Execute("13+3=", $true, $false, $false, $false, $false, $true, 1, $false, "25+3=", 2) | Out-Null
$d.Content.Find.Execute("3+3=", $true, $false, $false, $false, $false, $true, 1, $false, "30-13=", 2) | Out-Null
$d.Content.Find.Execute("96+2=", $true, $false, $false, $false, $false, $true, 1, $false, "2+82=", 2) | Out-Null
$d.Content.Find.Execute("1+41=", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=", 2) | Out-Null
$d.Content.Find.Execute("92-66=", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=", 2) | Out-Null
$d.Content.Find.Execute("73-12=", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=", 2) | Out-Null
$d.Content.Find.Execute("68-28=", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=", 2) | Out-Null
$d.Content.Find.Execute("16+4=", $true, $false, $false, $false, $false, $true, 1, $false, "75-16=", 2) | Out-Null
$d.Content.Find.Execute("89-32=", $true, $false, $false, $false, $false, $true, 1, $false, "97-0=", 2) | Out-Null
$d.Content.Find.Execute("76-73=", $true, $false, $false, $false, $false, $true, 1, $false, "35+55=", 2) | Out-Null
$d.Content.Find.Execute("98-74=", $true, $false, $false, $false, $false, $true, 1, $false, "66+18=", 2) | Out-Null
$d.Content.Find.Execute("49+11=", $true, $false, $false, $false, $false, $true, 1, $false, "86-24=", 2) | Out-Null
$d.Content.Find.Execute("75-36=", $true, $false, $false, $false, $false, $true, 1, $false, "30+52=", 2) | Out-Null
$d.Content.Find.Execute("76+18=", $true, $false, $false, $false, $false, $true, 1, $false, "86-17=", 2) | Out-Null
$d.Content.Find.Execute("4+86=", $true, $false, $false, $false, $false, $true, 1, $false, "36-6=", 2) | Out-Null
$d.Content.Find.Execute("52-13=", $true, $false, $false, $false, $false, $true, 1, $false, "93-89=", 2) | Out-Null
$d.Content.Find.Execute("1+30=", $true, $false, $false, $false, $false, $true, 1, $false, "42+26=", 2) | Out-Null
$d.Content.Find.Execute("38+35=", $true, $false, $false, $false, $false, $true, 1, $false, "13+7=", 2) | Out-Null
$d.Content.Find.Execute("15+6=", $true, $false, $false, $false, $false, $true, 1, $false, "38+35=", 2) | Out-Null
